# feat: add 2022-Q4 data
#
# 1) Insert a brand-new "2022-Q4" worksheet right after "总计" (i.e. right
#    before the existing "2022-Q3" sheet), formatted/laid out the same way
#    as the other quarterly fund-holding sheets, and fill it with the
#    Q4 figures.
# 2) Insert a new "2022-Q4" row at the top of the "总计" (totals) sheet's
#    data, pushing the existing Q3/Q2/Q1 rows down by one and renumbering
#    the index column.
#
# The "2022-Q3", "2022-Q2" and "2022-Q1" sheets themselves keep their
# names and data completely unchanged.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Build the new "2022-Q4" sheet, inserted right before "2022-Q3"
#    (i.e. right after "总计").
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (same headers/order as the other quarterly sheets)
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: index, code, name, scale, total position, position ratio,
# held market value (亿元), position rank.
$q4Rows = @(
    @(0, "290002", "泰信先行策略混合",     "5.69", "79.40", "6.84", "0.3892", 2),
    @(1, "013757", "泰信均衡价值混合A",     "0.67", "79.29", "8.49", "0.0569", 1),
    @(2, "005161", "华商上游产业股票",       "0.55", "88.76", "3.31", "0.0182", 7),
    @(3, "013758", "泰信均衡价值混合C",     "0.10", "79.29", "8.49", "0.0085", 1)
)

$r = 2
foreach ($row in $q4Rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    # Fund codes (column B) and columns D:G are stored as text (not
    # numeric) in this workbook, so force them to text with a leading
    # apostrophe to avoid auto-conversion to numbers (which would also
    # drop leading zeros from the fund codes).
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

$aCol = $q4.Range("A2:A5")
$aCol.Font.Bold = $true
$aCol.HorizontalAlignment = -4108
$aCol.VerticalAlignment = -4160
$aCol.Borders.LineStyle = 1

$q4.Range("A1").Select()

# ---------------------------------------------------------------------
# 2) Insert a new top row in "总计" for 2022-Q4, shifting the rest down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:A2").EntireRow.Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.47

$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

# Renumber the index column (A) for the rows that got pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

$total.Range("A1").Select()
